# Work diary and Documentation
# Append three new entries (rows 94-96) to the "JdT-TPI_LRD" journal sheet,
# matching the Tableau1 data table, and add a running-total formula in H96.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 94 ---------------------------------------------------------------
$ws.Range("A94").Value = 44711
$ws.Range("A94").NumberFormat = "dd/mm/yyyy"
$ws.Range("B94").Value = "Réalisation"
$ws.Range("C94").Value = 0.5
$ws.Range("D94").Value = "Corrections de bugs vus que pendant la mise en service"

# --- Row 95 ---------------------------------------------------------------
$ws.Range("A95").Value = 44711
$ws.Range("A95").NumberFormat = "dd/mm/yyyy"
$ws.Range("B95").Value = "Réalisation"
$ws.Range("C95").Value = 1
$ws.Range("D95").Value = "Constatations de nouveaux bugs et documentation de ses derniers"

# --- Row 96 ---------------------------------------------------------------
$ws.Range("A96").Value = 44711
$ws.Range("A96").NumberFormat = "dd/mm/yyyy"
$ws.Range("B96").Value = "Réalisation"
$ws.Range("C96").Value = 0.75
$ws.Range("D96").Value = "Documentation de la mise en service"

# Running total of hours worked, placed next to the last new row.
$ws.Range("H96").Formula = "=SUM(C:C)"

# Extend the "Tableau1" table (and its AutoFilter) to cover the new rows.
$lo = $ws.ListObjects.Item("Tableau1")
$lo.Resize($ws.Range("A1:F96"))

# Match the cursor position left behind by the author after the last edit.
$ws.Range("H97").Select()

$wb.Save()
